$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.65
$ws.Range("A6").Value = -22.056
$ws.Range("A7").Value = -21.286
$ws.Range("B7").Value = 6.312
$ws.Range("B12").Value = 5.065
$ws.Range("C13").Value = -13.225
$ws.Range("C14").Value = -11.992
$ws.Range("B15").Value = 5.087000000000001
$ws.Range("A16").Value = -21.771
$ws.Range("C16").Value = -13.14
$ws.Range("C19").Value = -12.228
$ws.Range("A20").Value = -21.433
$ws.Range("B20").Value = 6.121
$ws.Range("B21").Value = 8.904
$ws.Range("B22").Value = 8.516000000000002
$ws.Range("C22").Value = -12.393
$ws.Range("B23").Value = 7.359999999999999
$ws.Range("A28").Value = -21.687
$ws.Range("A29").Value = -21.675
$ws.Range("B29").Value = 5.731
$ws.Range("A32").Value = -21.685
$ws.Range("B34").Value = 8.018000000000001
$ws.Range("C36").Value = -12.37
$ws.Range("A40").Value = -20.452
$ws.Range("B42").Value = 7.582999999999998
$ws.Range("B43").Value = 5.671
$ws.Range("B44").Value = 4.993
$ws.Range("B45").Value = 5.048
$ws.Range("A46").Value = -21.15
$ws.Range("B46").Value = 6.344999999999999
$ws.Range("C46").Value = -13.898
$ws.Range("B50").Value = 5.291
$ws.Range("C50").Value = -13.363
$ws.Range("A51").Value = -20.943
$ws.Range("B51").Value = 7.513999999999998
$ws.Range("A52").Value = -21.262
$ws.Range("A57").Value = -22.263
$ws.Range("A59").Value = -22.415
$ws.Range("A62").Value = -21.782
$ws.Range("A66").Value = -21.551
$ws.Range("B66").Value = 5.473
$ws.Range("B67").Value = 5.249000000000001
$ws.Range("A73").Value = -20.426
$ws.Range("A74").Value = -21.045
$ws.Range("B79").Value = 5.53
$ws.Range("B84").Value = 5.798
$ws.Range("A92").Value = -21.33300000000001
$ws.Range("B92").Value = 5.527
$ws.Range("C95").Value = -11.719
$ws.Range("B97").Value = 6.811
$ws.Range("C97").Value = -13.606
$ws.Range("A100").Value = -21.481
